# Gorcery_Loads.xlsx update
# - Fix "Glenerin" -> "Glen Erin"
# - Insert new store rows (with their weekly delivery-time data left blank
#   for now, matching freshly added rows) at the correct positions so the
#   store list stays alphabetically/operationally grouped the way the
#   author laid it out
# - Append two new stores (Gateway / Voila) at the bottom with store
#   numbers 750 / 751
# - Update the current selection to reflect where the author left off
#   editing (O32)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing store name -------------------------------------------------
$ws.Range("B8").Value = "Glen Erin"

# --- insert new rows (top to bottom so row numbers below stay in sync) ------
# Each Insert() pushes the target row (and everything below it) down by one,
# and Excel carries the formatting of the row above into the newly blank row
# - which is exactly the formatting pattern the new rows need here.

# Before "BoxGrove" (row 15): Brookfield, Burloak
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()
$ws.Cells.Item(15, 1).Value = 15
$ws.Cells.Item(15, 2).Value = "Brookfield"
$ws.Cells.Item(16, 1).Value = 16
$ws.Cells.Item(16, 2).Value = "Burloak"

# After "BoxGrove" (now row 17), before "Aurora": Elizabeth St.
$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18, 1).Value = 18
$ws.Cells.Item(18, 2).Value = "Elizabeth St."

# After "Bathurst" (row 20), before "Richmond Hill": Bloor (HBC)
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).Value = 21
$ws.Cells.Item(21, 2).Value = "Bloor (HBC)"

# After "Richmond Hill" (row 22), before "MLS": First CDN Place
$ws.Rows.Item(23).Insert()
$ws.Cells.Item(23, 1).Value = 23
$ws.Cells.Item(23, 2).Value = "First CDN Place"

# After "Applewood", before "Guelph": Elizabeth, Mount Pleasant, Ancaster
$ws.Rows.Item(29).Insert()
$ws.Cells.Item(29, 1).Value = 29
$ws.Cells.Item(29, 2).Value = "Elizabeth"

$ws.Rows.Item(30).Insert()
$ws.Cells.Item(30, 1).Value = 30
$ws.Cells.Item(30, 2).Value = "Mount Pleasant"

$ws.Rows.Item(31).Insert()
$ws.Cells.Item(31, 1).Value = 31
$ws.Cells.Item(31, 2).Value = "Ancaster"

# After "Green Lane", before "Brooklin": Liberty Village
$ws.Rows.Item(37).Insert()
$ws.Cells.Item(37, 1).Value = 37
$ws.Cells.Item(37, 2).Value = "Liberty Village"

# After "Brooklin", before "Queensway": Bolton, Kitchener, Meadowvale, Kleinburg
$ws.Rows.Item(39).Insert()
$ws.Cells.Item(39, 1).Value = 39
$ws.Cells.Item(39, 2).Value = "Bolton"

$ws.Rows.Item(40).Insert()
$ws.Cells.Item(40, 1).Value = 40
$ws.Cells.Item(40, 2).Value = "Kitchener"

$ws.Rows.Item(41).Insert()
$ws.Cells.Item(41, 1).Value = 41
$ws.Cells.Item(41, 2).Value = "Meadowvale"

$ws.Rows.Item(42).Insert()
$ws.Cells.Item(42, 1).Value = 42
$ws.Cells.Item(42, 2).Value = "Kleinburg"

# After "Queensway" (now row 43): Colossus, then append Gateway / Voila
$ws.Rows.Item(44).Insert()
$ws.Cells.Item(44, 1).Value = 44
$ws.Cells.Item(44, 2).Value = "Colossus"

$ws.Cells.Item(45, 1).Value = 750
$ws.Cells.Item(45, 2).Value = "Gateway"
$ws.Cells.Item(46, 1).Value = 751
$ws.Cells.Item(46, 2).Value = "Voila"

# --- restore the author's last selection ------------------------------------
$ws.Range("O32").Select()
